# Auto-generated edit script applying cryptos.xlsx row updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.723.08"
$ws.Range("E2").Value = "  +0.56%  "

# Row 3
$ws.Range("D3").Value = "3.242.49"
$ws.Range("E3").Value = "  +1.37%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.91"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.06"
$ws.Range("E6").Value = "  +1.66%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "3.243.04"
$ws.Range("E8").Value = "  +1.31%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.547"
$ws.Range("E9").Value = "  +2.15%  "

# Row 10
$ws.Range("E10").Value = "  +0.18%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.67"
$ws.Range("E11").Value = "  -7.18%  "

# Row 12
$ws.Range("E12").Value = "  -0.48%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("E13").Value = "  +0.45%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.08"
$ws.Range("E14").Value = "  +0.07%  "

# Row 15
$ws.Range("D15").Value = "3.768.97"
$ws.Range("E15").Value = "  +1.41%  "

# Row 16
$ws.Range("D16").Value = "66.799.98"
$ws.Range("E16").Value = "  +0.86%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.45"
$ws.Range("E17").Value = "  +0.20%  "

# Row 18
$ws.Range("D18").Value = "3.233.89"
$ws.Range("E18").Value = "  +1.12%  "

# Row 19
$ws.Range("E19").Value = "  +1.16%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "513.43"
$ws.Range("E20").Value = "  +0.57%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.41"
$ws.Range("E21").Value = "  +0.63%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.739"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.11"
$ws.Range("E23").Value = "  +1.27%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.87"
$ws.Range("E24").Value = "  -1.75%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.84"
$ws.Range("E25").Value = "  -0.04%  "

# Row 26
$ws.Range("E26").Value = "  +0.10%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.49"
$ws.Range("E27").Value = "  +2.24%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.01"
$ws.Range("E28").Value = "  +0.77%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.41"
$ws.Range("E29").Value = "  +5.54%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.05"
$ws.Range("E30").Value = "  +4.78%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.13"
$ws.Range("E31").Value = "  +2.66%  "

# Row 32
$ws.Range("E32").Value = "  +0.16%  "

# Row 33
$ws.Range("E33").Value = "  +0.32%  "

# Row 34
$ws.Range("E34").Value = "  -3.24%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.55"
$ws.Range("E35").Value = "  -0.01%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "523.68"
$ws.Range("E36").Value = "  +8.01%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.50"
$ws.Range("E37").Value = "  +3.10%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0928"
$ws.Range("E38").Value = "  +2.97%  "

# Row 39
$ws.Range("D39").Value = "0.0₃0762"
$ws.Range("E39").Value = "  +16.45%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0420"
$ws.Range("E40").Value = "  +0.61%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.128"
$ws.Range("E41").Value = "  +4.69%  "

# Row 42
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.99"
$ws.Range("E42").Value = "  +2.42%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.82"
$ws.Range("E43").Value = "  -0.32%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.303"
$ws.Range("E44").Value = "  +1.86%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.52"
$ws.Range("E45").Value = "  +4.36%  "

# Row 46
$ws.Range("D46").Value = "2.866.30"
$ws.Range("E46").Value = "  -2.06%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.56"
$ws.Range("E47").Value = "  +0.35%  "

# Row 48
$ws.Range("E48").Value = "  +4.37%  "

# Row 49
$ws.Range("E49").Value = "  -0.05%  "

# Row 50
$ws.Range("E50").Value = "  +0.47%  "

# Row 51
$ws.Range("E51").Value = "  +1.74%  "
